$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.896.51'
$ws.Cells.Item(2, 5).Value = '  -3.11%  '
$ws.Cells.Item(3, 4).Value = '2.561.00'
$ws.Cells.Item(3, 5).Value = '  -1.53%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '505.61'
$ws.Cells.Item(5, 5).Value = '  -3.47%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '141.96'
$ws.Cells.Item(6, 5).Value = '  -7.86%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.553'
$ws.Cells.Item(8, 5).Value = '  -5.61%  '
$ws.Cells.Item(9, 4).Value = '2.565.56'
$ws.Cells.Item(9, 5).Value = '  -1.49%  '
$ws.Cells.Item(10, 5).Value = '  -7.28%  '
$ws.Cells.Item(11, 5).Value = '  -4.04%  '
$ws.Cells.Item(12, 5).Value = '  -4.83%  '
$ws.Cells.Item(13, 5).Value = '  -1.00%  '
$ws.Cells.Item(14, 4).Value = '3.006.86'
$ws.Cells.Item(14, 5).Value = '  -1.55%  '
$ws.Cells.Item(15, 4).Value = '58.901.21'
$ws.Cells.Item(15, 5).Value = '  -3.10%  '
$ws.Cells.Item(16, 5).Value = '  -5.03%  '
$ws.Cells.Item(17, 5).Value = '  -4.56%  '
$ws.Cells.Item(18, 4).Value = '2.582.80'
$ws.Cells.Item(18, 5).Value = '  -0.81%  '
$ws.Cells.Item(19, 5).Value = '  -5.40%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '331.68'
$ws.Cells.Item(20, 5).Value = '  -6.61%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '10.03'
$ws.Cells.Item(21, 5).Value = '  -4.92%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  -0.03%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.93'
$ws.Cells.Item(23, 5).Value = '  -4.36%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '59.50'
$ws.Cells.Item(24, 5).Value = '  -2.58%  '
$ws.Cells.Item(25, 5).Value = '  -4.97%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.999'
$ws.Cells.Item(26, 5).Value = '  -0.02%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.158'
$ws.Cells.Item(27, 5).Value = '  -4.93%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0775'
$ws.Cells.Item(28, 5).Value = '  -7.89%  '
$ws.Cells.Item(29, 5).Value = '  -7.20%  '
$ws.Cells.Item(30, 5).Value = '  -0.04%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '148.99'
$ws.Cells.Item(31, 5).Value = '  -0.43%  '
$ws.Cells.Item(32, 5).Value = '  -4.74%  '
$ws.Cells.Item(33, 5).Value = '  -4.11%  '
$ws.Cells.Item(34, 5).Value = '  -7.87%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.86'
$ws.Cells.Item(35, 5).Value = '  -7.79%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.875'
$ws.Cells.Item(36, 5).Value = '  -4.81%  '
$ws.Cells.Item(37, 5).Value = '  -8.11%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '35.85'
$ws.Cells.Item(38, 5).Value = '  -1.64%  '
$ws.Cells.Item(39, 5).Value = '  -9.57%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '286.23'
$ws.Cells.Item(40, 5).Value = '  -3.88%  '
$ws.Cells.Item(41, 5).Value = '  -7.55%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.49'
$ws.Cells.Item(42, 5).Value = '  -7.63%  '
$ws.Cells.Item(43, 5).Value = '  +0.03%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0980'
$ws.Cells.Item(44, 5).Value = '  -3.25%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.606'
$ws.Cells.Item(45, 5).Value = '  -2.65%  '
$ws.Cells.Item(46, 5).Value = '  -5.51%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.34'
$ws.Cells.Item(47, 5).Value = '  +0.09%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '18.59'
$ws.Cells.Item(48, 5).Value = '  -4.99%  '
$ws.Cells.Item(49, 5).Value = '  -5.37%  '
$ws.Cells.Item(50, 5).Value = '  -8.07%  '
$ws.Cells.Item(51, 4).Value = '1.882.16'
$ws.Cells.Item(51, 5).Value = '  -4.18%  '
